$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new analysis columns: F = "Årsag" (Reason), G = "Ny leverandør" (New supplier) ---
# The existing F column ("TCV_range") is moved to the new H column, and the old F values
# (which were all "180000-200000") are replaced with per-row "Årsag" answers.

# 1) Copy header formatting from the existing F1 header cell into the two new header cells
#    (G1 and H1) so they match the other bold/bordered header cells.
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("F1").Copy($ws.Range("H1"))

# 2) Move the old "TCV_range" value column (currently F2:F5) down into the new H2:H5 cells.
$ws.Range("H2").Value = $ws.Cells.Item(2, 6).Value2
$ws.Range("H3").Value = $ws.Cells.Item(3, 6).Value2
$ws.Range("H4").Value = $ws.Cells.Item(4, 6).Value2
$ws.Range("H5").Value = $ws.Cells.Item(5, 6).Value2

# 3) Re-purpose the F column header to "Årsag" and give each row its reason value.
$ws.Range("F1").Value = "Årsag"
$ws.Range("F2").Value = "Strategisk beslutning"
$ws.Range("F3").Value = "Fusionerer med anden virksomhed"
$ws.Range("F4").Value = "Utilfredshed (Service - uddyb i bemærkninger)"
$ws.Range("F5").Value = "Pris"

# 4) New G column header - "Ny leverandør" (New supplier); data cells stay empty.
$ws.Range("G1").Value = "Ny leverandør"

# 5) New H column header - "TCV_range" (values were already populated above).
$ws.Range("H1").Value = "TCV_range"
